# Update the "取得日時" (acquisition timestamp) column (A) for rows 2-11
# on the first worksheet ("ランサーズ") from the old timestamp
# "2025-09-14 01:19:48" to the new timestamp "2025-09-14 01:48:20".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-09-14 01:19:48"
$newValue = "2025-09-14 01:48:20"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
